$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.255.50"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "2.643.00"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.90"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.64"
$ws.Range("E6").Value = "  +3.38%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.142"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.25"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.350"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.87"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "3.126.51"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("D16").Value = "68.197.30"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").Value = "2.650.02"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.36"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "360.14"
$ws.Range("E19").Value = "  -1.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.39"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.73"
$ws.Range("E22").Value = "  -3.08%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.24"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.79"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").Value = "2.777.59"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000103"
$ws.Range("E28").Value = "  -3.25%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "558.39"
$ws.Range("E30").Value = "  -3.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.95"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.37"
$ws.Range("E32").Value = "  -3.80%  "
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -3.29%  "
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.71"
$ws.Range("E37").Value = "  +1.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.74"
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.85"
$ws.Range("E40").Value = "  -2.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.31"
$ws.Range("E41").Value = "  -2.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.61"
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("D43").Value = "0.0₆0321"
$ws.Range("E43").Value = "  -7.87%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "156.75"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.75"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.95"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0773"
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.564"
$ws.Range("E51").Value = "  -0.44%  "
